# "Generate Report for handoff"
#
# The previous handoff attempt (source file
# 51853478-e905-4c66-8e2c-0f054f664e21.md) is replaced by a new report row
# for cbae07d2-d130-47c4-838c-a3cc4cb73d0d.md whose handoff transform
# failed. The per-language sheets (zh-cn, de-de) lose their "Latest Handoff
# File" link/value (handoff never produced a target file), their handoff
# datetime + handback datetime reset to the zero date, and the handoff
# reason flips from "Include" to "Ignored". The Overview sheet mirrors the
# new file name + failed status for both languages.

$wb = $excel.ActiveWorkbook

$oldName = "51853478-e905-4c66-8e2c-0f054f664e21.md"
$newName = "cbae07d2-d130-47c4-838c-a3cc4cb73d0d.md"
$oldStatus = "Ready for handoff"
$newStatus = "Handoff transform failed"
$zeroDate = "0001-01-01 00:00:00"
$oldReason = "Include"
$newReason = "Ignored"

# ---------------------------------------------------------------------
# Overview sheet: update the file-name hyperlink and the two status cells
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value2 = $newName
foreach ($hl in $ov.Hyperlinks) {
    if ($hl.Range.Address() -eq '$A$2') {
        $hl.TextToDisplay = $newName
        $hl.Address = $hl.Address.Replace($oldName, $newName)
    }
}

$ov.Range("B2").Value2 = $newStatus
$ov.Range("C2").Value2 = $newStatus

# ---------------------------------------------------------------------
# Per-language sheets: zh-cn, de-de
# ---------------------------------------------------------------------
$langSheets = @("zh-cn", "de-de")

foreach ($sheetName in $langSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    # A2: file name + hyperlink
    $ws.Range("A2").Value2 = $newName
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq '$A$2') {
            $hl.TextToDisplay = $newName
            $hl.Address = $hl.Address.Replace($oldName, $newName)
        }
    }

    # B2: status
    $ws.Range("B2").Value2 = $newStatus

    # C2: Latest Handoff File - transform failed, so no handoff file anymore.
    # Remove the hyperlink and clear the cell (and its hyperlink styling).
    $toDelete = @()
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq '$C$2') {
            $toDelete += $hl
        }
    }
    foreach ($hl in $toDelete) {
        $hl.Delete()
    }
    $ws.Range("C2").ClearContents()
    $ws.Range("C2").Style = "Normal"

    # D2: Latest Handoff Datetime - reset to zero date (handoff never completed)
    $ws.Range("D2").Value2 = $zeroDate

    # G2: Latest Handback DateTime - stays the zero date (unchanged value,
    # rewritten for parity with the regenerated report)
    $ws.Range("G2").Value2 = $zeroDate

    # H2: Handoff Reason - Include -> Ignored
    $ws.Range("H2").Value2 = $newReason

    # Row 3 (.localization-config) is informational only and carries no
    # semantic change in this report refresh, aside from re-asserting its
    # already-correct zero-date/Ignored values.
    $ws.Range("D3").Value2 = $zeroDate
    $ws.Range("G3").Value2 = $zeroDate
    $ws.Range("H3").Value2 = $newReason
}
